$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "depositType"/"depositValue"/"profitType"/"profitValue" columns (K:N)
# are being removed from the task table, along with the header cell comment
# on M1 ("profitType") that documents them. Remove that comment first so it
# does not linger orphaned once its column is gone.
$ws.Range("M1").Comment.Delete()

# A couple of existing rows get their taskStyleId bumped (priority of a
# "send letter" style task), and the priority column values (J) are updated
# for several rows.
$ws.Range("C7").Value = 2
$ws.Range("C8").Value = 2

$ws.Range("J3").Value = 19
$ws.Range("J4").Value = 100
$ws.Range("J5").Value = 10
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 100
$ws.Range("J8").Value = 100

# Remove the now-unused depositType/depositValue/profitType/profitValue
# columns (K:N), shifting everything after them to the left.
$ws.Range("K1:N16").EntireColumn.Delete()

# Mirror the author's final selection state (whole column J selected).
$null = $ws.Columns.Item(10).Select()
